# Refresh the crypto price/volume snapshot for Sheet1 (GitHub Actions scheduled
# update, commit "Updated cryptos list on Thu Mar 28 16:28:13 UTC 2024").
#
# Column layout: A=rank(unchanged) B=Coin C=Link D=Price E=Volume(1h)
#
# Two coin-ranking swaps happened between snapshots (their relative market
# positions crossed), so rows 17/18 and 43/44 get their whole B:E payload
# exchanged (not just the price/volume refreshed in place).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Many "Price" cells are plain decimal-looking strings (e.g. "582.84",
# "0.0450"). Assigning those to .Value directly lets Excel auto-coerce them
# into real numbers (and "0.0450" would lose its trailing zero as 0.045), so
# for those we force the Text number format before writing, then drop back to
# the built-in "Normal" cell style so no stray formatting is left behind.
function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Row 2
$ws.Range("D2").Value = "70.797.45"
$ws.Range("E2").Value = "  +2.56%  "

# Row 3
$ws.Range("D3").Value = "3.566.07"
$ws.Range("E3").Value = "  +1.71%  "

# Row 4
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
Set-TextValue "D5" "582.84"
$ws.Range("E5").Value = "  +2.18%  "

# Row 6
Set-TextValue "D6" "186.24"
$ws.Range("E6").Value = "  +2.42%  "

# Row 7
$ws.Range("D7").Value = "3.553.82"
$ws.Range("E7").Value = "  +1.52%  "

# Row 8
$ws.Range("E8").Value = "  +1.42%  "

# Row 9
$ws.Range("E9").Value = "  -0.08%  "

# Row 10
Set-TextValue "D10" "0.220"
$ws.Range("E10").Value = "  +16.96%  "

# Row 11
$ws.Range("E11").Value = "  +2.58%  "

# Row 12
Set-TextValue "D12" "54.67"
$ws.Range("E12").Value = "  +1.85%  "

# Row 13
Set-TextValue "D13" "0.0000319"
$ws.Range("E13").Value = "  +6.03%  "

# Row 14
Set-TextValue "D14" "9.53"
$ws.Range("E14").Value = "  +1.26%  "

# Row 15
$ws.Range("D15").Value = "4.033.67"
$ws.Range("E15").Value = "  -0.86%  "

# Row 16
$ws.Range("D16").Value = "70.790.33"
$ws.Range("E16").Value = "  +2.96%  "

# Row 17
$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue "D17" "19.33"
$ws.Range("E17").Value = "  +0.71%  "

# Row 18
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.547.13"
$ws.Range("E18").Value = "  +1.48%  "

# Row 19
Set-TextValue "D19" "12.43"
$ws.Range("E19").Value = "  -0.11%  "

# Row 20
Set-TextValue "D20" "565.26"
$ws.Range("E20").Value = "  +6.08%  "

# Row 21
$ws.Range("E21").Value = "  +0.87%  "

# Row 22
$ws.Range("E22").Value = "  -1.64%  "

# Row 23
Set-TextValue "D23" "17.80"
$ws.Range("E23").Value = "  -11.69%  "

# Row 24
Set-TextValue "D24" "5.04"
$ws.Range("E24").Value = "  +1.25%  "

# Row 25
Set-TextValue "D25" "4.58"
$ws.Range("E25").Value = "  +5.20%  "

# Row 26
Set-TextValue "D26" "94.15"
$ws.Range("E26").Value = "  +0.41%  "

# Row 27
Set-TextValue "D27" "11.31"
$ws.Range("E27").Value = "  +3.05%  "

# Row 28
Set-TextValue "D28" "2.95"
$ws.Range("E28").Value = "  +2.04%  "

# Row 29
Set-TextValue "D29" "9.19"
$ws.Range("E29").Value = "  +0.94%  "

# Row 30
Set-TextValue "D30" "32.50"
$ws.Range("E30").Value = "  +3.03%  "

# Row 31
Set-TextValue "D31" "7.29"
$ws.Range("E31").Value = "  +0.33%  "

# Row 32
Set-TextValue "D32" "12.35"
$ws.Range("E32").Value = "  -1.78%  "

# Row 33
Set-TextValue "D33" "0.117"
$ws.Range("E33").Value = "  +2.90%  "

# Row 34
Set-TextValue "D34" "63.61"
$ws.Range("E34").Value = "  -0.61%  "

# Row 35
Set-TextValue "D35" "3.37"
$ws.Range("E35").Value = "  +8.52%  "

# Row 36
Set-TextValue "D36" "553.77"
$ws.Range("E36").Value = "  -2.55%  "

# Row 37
Set-TextValue "D37" "0.419"
$ws.Range("E37").Value = "  +5.49%  "

# Row 38
$ws.Range("D38").Value = "0.0₃0803"
$ws.Range("E38").Value = "  +5.58%  "

# Row 39
Set-TextValue "D39" "37.77"
$ws.Range("E39").Value = "  -0.45%  "

# Row 40
$ws.Range("E40").Value = "  -0.05%  "

# Row 41
Set-TextValue "D41" "3.37"
$ws.Range("E41").Value = "  +10.30%  "

# Row 42
$ws.Range("D42").Value = "3.550.13"

# Row 43
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D43" "3.45"
$ws.Range("E43").Value = "  +3.28%  "

# Row 44
$ws.Range("B44").Value = "Kaspa"
$ws.Range("C44").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D44" "0.137"
$ws.Range("E44").Value = "  +3.24%  "

# Row 45
Set-TextValue "D45" "0.0450"

# Row 46
Set-TextValue "D46" "3.49"
$ws.Range("E46").Value = "  +0.31%  "

# Row 47
$ws.Range("E47").Value = "  -0.94%  "

# Row 48
Set-TextValue "D48" "9.36"
$ws.Range("E48").Value = "  +2.18%  "

# Row 49
$ws.Range("E49").Value = "  +2.77%  "

# Row 50
Set-TextValue "D50" "1.48"
$ws.Range("E50").Value = "  +9.06%  "

# Row 51
Set-TextValue "D51" "0.999"
$ws.Range("E51").Value = "  +0.21%  "
